$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the existing answer text in B2:B4 (old "A"/"B"/"Telkom") with the
#    new, longer answer text.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Artificial Intelligence has a wide-ranging impact on various aspects of our lives"
$ws.Range("B3").Value = "Artificial Intelligence is perceived as including elements of imagination and it is based on futuristic ideas/the threat that AI might take over the world."
$ws.Range("B4").Value = "Artificial intelligence is a discipline in computer science in which machines copy the way humans think and make decisions."

# ---------------------------------------------------------------------------
# 2. Rows 5 and 6 (old "Marking App"/"Geek") both become the single letter "A".
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "A"
$ws.Range("B6").Value = "A"

# ---------------------------------------------------------------------------
# 3. Append four brand new rows (7-10), copying the question-number / answer
#    border formatting from the row directly above so the new rows look the
#    same as the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A7").Interior.ColorIndex = -4142   # xlNone
[void]$ws.Range("A6:B6").Copy()

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Pretoria"

$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Interior.ColorIndex = -4142
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "To highlight/show the rapid growth of Artificial Intelligence and its importance in businesses."

$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Interior.ColorIndex = -4142
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "The words suggest that using Artificial Intelligence brings about amazing reforms/innovations."

$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Interior.ColorIndex = -4142
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "South Africa"

# ---------------------------------------------------------------------------
# 4. Widen column B so the long answer text fits, and update the selected
#    cell to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 123

[void]$ws.Range("C7").Select()
